$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a scratch cell far outside the used range to stage the literal text
# "2008-05-17" as a formula result, then paste-special (values only) it into
# each BF cell. This avoids Excel's automatic "looks like a date" type
# coercion that a direct .Value assignment would trigger, and it does not
# touch the destination cell's number format / style (PasteSpecial values-only
# only moves the value, not formatting).
$scratch = $ws.Cells.Item(1, 100)
$scratch.Formula = "=""2008-05-17"""
$scratch.Copy()

for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 58).PasteSpecial(-4163)
}

$excel.CutCopyMode = $false
$scratch.Clear()
